$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was inserted as the first data row of the
# "Poroto granado" time series (row 30), pushing every following row
# (old 30..113) down by one (new 31..114).
$ws.Rows.Item(30).Insert()

$ws.Range("A30").Value = 8
$ws.Range("B30").Value = "Terminal La Palmera de La Serena"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44953
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 100112030
$ws.Range("G30").Value = "Poroto granado"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 440
$ws.Range("K30").Value = 37000
$ws.Range("L30").Value = 38000
$ws.Range("M30").Value = 37500
$ws.Range("N30").Value = "$/malla 25 kilos"
$ws.Range("O30").Value = "Provincia del Elquí"
$ws.Range("P30").Value = 1500
$ws.Range("Q30").Value = 25
$ws.Range("R30").Value = "Hortaliza"
